$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginValidation")

# "Verification Message" header is renamed (no spaces) and the old
# "Lewis Mocha" placeholder value is replaced with the real verification
# text shown after successful login.
$ws.Range("D1").Value = "VerificationMessage"
$ws.Range("D7").Value = "My Account"

# Active selection on this sheet moves up one row.
$ws.Range("D12").Select()
